# Applies the "updated slides and literature" commit:
#   - TextBox 17 (slide 2): grows taller (new bottom bullets) and shifts down
#   - Straight Arrow Connector 42 (slide 2): gains a touch of slope (width/height)
#   - Group 49 (slide 2): nudges down-right
#
# NOTE: Left/Top/Width/Height on Shape are expressed in points (1 pt = 12700 EMU)
# and are rounded through a single-precision float internally, so the literal
# point values below are chosen (via EMU/12700, fine-tuned by float32 ULPs) so
# that they land exactly on the target EMU values after conversion.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- TextBox 17: reposition/resize, then append two new bullet paragraphs ---
$textBox = $s.Shapes.Item(2)
$textBox.Top = 182.95008850097656
$textBox.Height = 109.05472564697266

$tr = $textBox.TextFrame.TextRange
$lastBullet = $tr.Paragraphs(4, 1)
$lastBullet.InsertAfter("`rImplicit neural representation.") | Out-Null
$newBullet = $tr.Paragraphs(5, 1)
$newBullet.InsertAfter("`rDiffusion models.") | Out-Null

# --- Straight Arrow Connector 42: widen and add a bit of drop ---
$connector = $s.Shapes.Item(7)
$connector.Width = 114.67024230957031
$connector.Height = 2.6412599086761475

# --- Group 49: shift down and to the right ---
$group = $s.Shapes.Item(12)
$group.Left = 716.8335571289062
$group.Top = 20.222599029541016
